# Updates the scale_hx711 BoM sheet to match the regenerated KiCad BOM export:
#  - refresh the generation timestamp and component count
#  - insert a new connector J4 (Conn_01x04) after J3
#  - insert a new transistor Q4 (AO3400) after Q3
#  - update R5's footprint and R6/R7's values
#  - insert two new resistors R14, R15 after R13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header metadata -------------------------------------------------
$ws.Range("B2").Value = "Tue Mar 10 15:40:08 2020"
$ws.Range("B4").Value = 79

# --- Insert J4 after J3 (row 47) -> new row 48 ------------------------
$ws.Rows.Item(48).Insert()
$ws.Range("A48").Value = "J4"
$ws.Range("B48").Value = "Conn_01x04"
$ws.Range("C48").Value = "Connector_PinHeader_2.00mm:PinHeader_1x04_P2.00mm_Vertical"
$ws.Range("D48").Value = "~"

# --- Insert Q4 after Q3 (now row 59) -> new row 60 --------------------
$ws.Rows.Item(60).Insert()
$ws.Range("A60").Value = "Q4"
$ws.Range("B60").Value = "AO3400"
$ws.Range("C60").Value = "Package_TO_SOT_SMD:SOT-23"

# --- Update R5 footprint, R6 and R7 values (now rows 65-67) -----------
# R6/R7 previously held plain numbers (style carried a left-aligned numeric
# format); resetting HorizontalAlignment back to "general" drops that
# inherited formatting now that the cells hold text again.
$ws.Range("C65").Value = "Resistor_SMD:R_1206_3216Metric"
$ws.Range("B66").HorizontalAlignment = 1
$ws.Range("B66").Value = "4.7k"
$ws.Range("B67").HorizontalAlignment = 1
$ws.Range("B67").Value = "24k"

# --- Insert R14, R15 after R13 (now row 73) -> new rows 74, 75 --------
$ws.Rows.Item(74).Insert()
$ws.Range("A74").Value = "R14"
$ws.Range("B74").Value = "10k"
$ws.Range("C74").Value = "Resistor_SMD:R_0603_1608Metric_Pad1.05x0.95mm_HandSolder"

$ws.Rows.Item(75).Insert()
$ws.Range("A75").Value = "R15"
$ws.Range("B75").Value = "10k"
$ws.Range("C75").Value = "Resistor_SMD:R_0603_1608Metric_Pad1.05x0.95mm_HandSolder"
